$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

# Columns that change identically for both data rows (2 and 3)
$common = @{
    "D"  = -0.153
    "E"  = -0.047
    "G"  = -0.2772277227722773
    "H"  = -0.2772277227722773
    "I"  = -0.4099009900990099
    "J"  = -0.4099009900990099
    "K"  = 3.25
    "L"  = 3.217821782178218
    "U"  = 2.89
    "V"  = 0.07727272727272728
    "W"  = 0.1212686567164179
    "X"  = 0.1128686424647537
    "Y"  = 0.008400014251664173
    "Z"  = 0.03787878787878788
    "AA" = -0.01552655265526552
    "AB" = 0.1123939582790026
    "AC" = -0.1279205109342682
    "AD" = 0.349
    "AF" = 0.349
    "AG" = -2.541
    "AH" = 0.009245278020609819
    "AI" = 0.01161436320676229
    "AJ" = -0.07289365730514361
    "AK" = -0.09356014580801945
    "AL" = 0.053
    "AM" = -0.01
    "AN" = -0.9614325068870523
    "AO" = -7.811320754716981
    "AP" = 7.000000000000001
    "AQ" = 41.39999999999999
}

foreach ($col in $common.Keys) {
    $ws.Range("$col`2").Value = $common[$col]
    $ws.Range("$col`3").Value = $common[$col]
}

# cash_returned block: row 2 settles at positive zero, row 3 at negative zero
$zeroCols = @("M", "N", "O", "P", "Q", "R")
foreach ($col in $zeroCols) {
    $ws.Range("$col`2").Value = 0
    $ws.Range("$col`3").Value = -0.0
}

# buybacks_cash_returned column (T) is removed entirely from the sheet
$ws.Range("T2").ClearContents()
$ws.Range("T3").ClearContents()
